$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.472.01'
$ws.Range('E2').Value = '  -1.98%  '
$ws.Range('D3').Value = '1.954.28'
$ws.Range('E3').Value = '  -0.50%  '
$ws.Range('D4').Value = '''1.012'
$ws.Range('E4').Value = '  +0.45%  '
$ws.Range('D5').Value = '''322.00'
$ws.Range('E5').Value = '  -1.47%  '
$ws.Range('D6').Value = '''1.010'
$ws.Range('E6').Value = '  +0.40%  '
$ws.Range('D7').Value = '''0.4791'
$ws.Range('E7').Value = '  -4.00%  '
$ws.Range('D8').Value = '''0.4077'
$ws.Range('E8').Value = '  -2.89%  '
$ws.Range('D9').Value = '''54.06'
$ws.Range('E9').Value = '  +3.41%  '
$ws.Range('D10').Value = '''0.08495'
$ws.Range('E10').Value = '  -7.23%  '
$ws.Range('D11').Value = '''1.058'
$ws.Range('E11').Value = '  -3.52%  '
$ws.Range('D12').Value = '''22.43'
$ws.Range('E12').Value = '  -1.96%  '
$ws.Range('D13').Value = '1.986.94'
$ws.Range('E13').Value = '  -1.44%  '
$ws.Range('D14').Value = '''7.571'
$ws.Range('E14').Value = '  -3.52%  '
$ws.Range('D15').Value = '''6.160'
$ws.Range('E15').Value = '  -4.21%  '
$ws.Range('E16').Value = '  +0.54%  '
$ws.Range('D17').Value = '''90.57'
$ws.Range('E17').Value = '  -0.78%  '
$ws.Range('D18').Value = '''0.00001072'
$ws.Range('E18').Value = '  -2.51%  '
$ws.Range('D19').Value = '''0.06635'
$ws.Range('E19').Value = '  -0.84%  '
$ws.Range('D20').Value = '''18.47'
$ws.Range('E20').Value = '  -3.68%  '
$ws.Range('D21').Value = '''1.009'
$ws.Range('E21').Value = '  +0.31%  '
$ws.Range('D22').Value = '''5.838'
$ws.Range('E22').Value = '  -2.15%  '
$ws.Range('D23').Value = '28.475.95'
$ws.Range('E23').Value = '  -1.96%  '
$ws.Range('D24').Value = '''11.44'
$ws.Range('E24').Value = '  -5.02%  '
$ws.Range('D25').Value = '''2.299'
$ws.Range('E25').Value = '  +0.70%  '
$ws.Range('D26').Value = '2.170.67'
$ws.Range('E26').Value = '  -3.00%  '
$ws.Range('D27').Value = '''156.25'
$ws.Range('E27').Value = '  -0.08%  '
$ws.Range('D28').Value = '''20.29'
$ws.Range('E28').Value = '  -1.30%  '
$ws.Range('D29').Value = '''2.175'
$ws.Range('E29').Value = '  -3.73%  '
$ws.Range('D30').Value = '''5.802'
$ws.Range('E30').Value = '  -5.81%  '
$ws.Range('D31').Value = '''124.36'
$ws.Range('E31').Value = '  -1.81%  '
$ws.Range('D32').Value = '''0.9855'
$ws.Range('E32').Value = '  -4.92%  '
$ws.Range('D33').Value = '''0.09656'
$ws.Range('E33').Value = '  -1.75%  '
$ws.Range('D34').Value = '''1.448'
$ws.Range('E34').Value = '  -4.87%  '
$ws.Range('D35').Value = '''3.693'
$ws.Range('E35').Value = '  +0.43%  '
$ws.Range('D36').Value = '''5.626'
$ws.Range('E36').Value = '  -2.25%  '
$ws.Range('D37').Value = '''9.122'
$ws.Range('E37').Value = '  +2.16%  '
$ws.Range('D38').Value = '''0.02336'
$ws.Range('E38').Value = '  -3.19%  '
$ws.Range('D39').Value = '''0.06192'
$ws.Range('E39').Value = '  -2.08%  '
$ws.Range('D40').Value = '''1.252'
$ws.Range('E40').Value = '  -3.92%  '
$ws.Range('D41').Value = '''0.6227'
$ws.Range('E41').Value = '  -3.18%  '
$ws.Range('E42').Value = '  -1.92%  '
$ws.Range('E43').Value = '  +0.36%  '
$ws.Range('D44').Value = '''0.1919'
$ws.Range('E44').Value = '  -2.95%  '
$ws.Range('D45').Value = '''1.327'
$ws.Range('E45').Value = '  +3.05%  '
$ws.Range('D46').Value = '''0.5951'
$ws.Range('E46').Value = '  -4.18%  '
$ws.Range('D47').Value = '''13.01'
$ws.Range('E47').Value = '  -2.27%  '
$ws.Range('D48').Value = '''2.058'
$ws.Range('E48').Value = '  -5.27%  '
$ws.Range('D49').Value = '''3.410'
$ws.Range('E49').Value = '  -1.46%  '
$ws.Range('D50').Value = '''0.06815'
$ws.Range('E50').Value = '  -1.92%  '
$ws.Range('B51').Value = 'Quant'
$ws.Range('C51').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D51').Value = '''111.09'
$ws.Range('E51').Value = '  -1.26%  '
